$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (rId1 / sheet1.xml) — 想去人数 (F column) updates
$ws1.Range("F2").Value = 3125
$ws1.Range("F3").Value = 525
$ws1.Range("F4").Value = 1086
$ws1.Range("F5").Value = 79
$ws1.Range("F6").Value = 29
$ws1.Range("F8").Value = 33
$ws1.Range("F9").Value = 1118
$ws1.Range("F10").Value = 15557
$ws1.Range("F12").Value = 166
$ws1.Range("F13").Value = 1020
$ws1.Range("F14").Value = 6138
$ws1.Range("F16").Value = 106
$ws1.Range("F19").Value = 109
$ws1.Range("F22").Value = 628
$ws1.Range("F24").Value = 8
$ws1.Range("F25").Value = 2
$ws1.Range("F26").Value = 209
$ws1.Range("F27").Value = 860
$ws1.Range("F28").Value = 22
$ws1.Range("F29").Value = 4992
$ws1.Range("F30").Value = 474
$ws1.Range("F31").Value = 11004
$ws1.Range("F33").Value = 10
$ws1.Range("F34").Value = 110
$ws1.Range("F35").Value = 157
$ws1.Range("F36").Value = 3793

# Sheet "全部类型" (rId4 / sheet4.xml) — 想去人数 (F column) updates
$ws4.Range("F3").Value = 3125
$ws4.Range("F4").Value = 525
$ws4.Range("F5").Value = 1086
$ws4.Range("F6").Value = 79
$ws4.Range("F7").Value = 29
$ws4.Range("F9").Value = 33
$ws4.Range("F10").Value = 1118
$ws4.Range("F11").Value = 15557
$ws4.Range("F13").Value = 166
$ws4.Range("F14").Value = 1020
$ws4.Range("F15").Value = 6138
$ws4.Range("F17").Value = 106
$ws4.Range("F20").Value = 109
$ws4.Range("F23").Value = 628
$ws4.Range("F25").Value = 8
$ws4.Range("F26").Value = 2
$ws4.Range("F27").Value = 209
$ws4.Range("F28").Value = 859
$ws4.Range("F29").Value = 22
$ws4.Range("F30").Value = 4992
$ws4.Range("F31").Value = 474
$ws4.Range("F33").Value = 11004
$ws4.Range("F35").Value = 10
$ws4.Range("F36").Value = 110
$ws4.Range("F37").Value = 157
$ws4.Range("F38").Value = 3793
